$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: fix date in A1 (advance one day)
$ws.Range("A1").Value = 45309

# Step 2: update prices in D29 / D30
$ws.Range("D29").Value = 364.992
$ws.Range("D30").Value = 514.29
